# Auto-generated COM script: applies the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.563.92'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '1.472.22'
$ws.Range("E3").Value = '  +2.24%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '0.9576'
$ws.Range("E5").Value = '  +3.84%  '
$ws.Range("D6").Value = '277.05'
$ws.Range("E6").Value = '  +0.86%  '
$ws.Range("D7").Value = '0.3549'
$ws.Range("E7").Value = '  -2.62%  '
$ws.Range("D8").Value = '0.3069'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '1.089'
$ws.Range("E9").Value = '  +7.31%  '
$ws.Range("D10").Value = '39.45'
$ws.Range("E10").Value = '  +1.80%  '
$ws.Range("D11").Value = '0.06632'
$ws.Range("E11").Value = '  +2.61%  '
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.10'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +4.61%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.453'
$ws.Range("E14").Value = '  +2.77%  '
$ws.Range("D15").Value = '6.172'
$ws.Range("E15").Value = '  +2.58%  '
$ws.Range("D16").Value = '0.9591'
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("D18").Value = '1.468.48'
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("D19").Value = '0.05963'
$ws.Range("E19").Value = '  +5.72%  '
$ws.Range("D20").Value = '68.77'
$ws.Range("E20").Value = '  +1.56%  '
$ws.Range("D21").Value = '5.477'
$ws.Range("E21").Value = '  +2.93%  '
$ws.Range("D22").Value = '14.47'
$ws.Range("E22").Value = '  +2.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.20'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +4.56%  '
$ws.Range("D24").Value = '2.279'
$ws.Range("E24").Value = '  +1.74%  '
$ws.Range("D25").Value = '20.589.89'
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("D26").Value = '145.59'
$ws.Range("E26").Value = '  +4.56%  '
$ws.Range("D27").Value = '2.086'
$ws.Range("E27").Value = '  +2.61%  '
$ws.Range("D28").Value = '17.14'
$ws.Range("E28").Value = '  +1.85%  '
$ws.Range("D29").Value = '1.627.59'
$ws.Range("E29").Value = '  +2.17%  '
$ws.Range("D30").Value = '114.35'
$ws.Range("E30").Value = '  +3.90%  '
$ws.Range("D31").Value = '3.983'
$ws.Range("E31").Value = '  -1.98%  '
$ws.Range("D32").Value = '4.929'
$ws.Range("E32").Value = '  +3.09%  '
$ws.Range("D33").Value = '0.07915'
$ws.Range("E33").Value = '  +3.73%  '
$ws.Range("D34").Value = '0.7927'
$ws.Range("E34").Value = '  +2.17%  '
$ws.Range("D35").Value = '1.207'
$ws.Range("E35").Value = '  +8.16%  '
$ws.Range("D36").Value = '1.438'
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05670'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("D38").Value = '4.711'
$ws.Range("E38").Value = '  +1.88%  '
$ws.Range("D39").Value = '0.9592'
$ws.Range("E39").Value = '  +2.46%  '
$ws.Range("D40").Value = '0.02025'
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("D41").Value = '10.28'
$ws.Range("D42").Value = '0.1847'
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("D43").Value = '7.303'
$ws.Range("E43").Value = '  +5.93%  '
$ws.Range("D44").Value = '3.513'
$ws.Range("E44").Value = '  +1.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5220'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("D46").Value = '12.03'
$ws.Range("E46").Value = '  +3.25%  '
$ws.Range("D47").Value = '120.16'
$ws.Range("E47").Value = '  +4.99%  '
$ws.Range("D48").Value = '0.5166'
$ws.Range("E48").Value = '  +2.53%  '
$ws.Range("D49").Value = '1.803'
$ws.Range("E49").Value = '  +4.95%  '
$ws.Range("D50").Value = '0.06409'
$ws.Range("E50").Value = '  +0.81%  '
$ws.Range("D51").Value = '0.9908'
$ws.Range("E51").Value = '  +0.87%  '
